$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header relabeling: shift "max" text from C1 to E1, and move
# "prediction"/"rejection-f" up into C1/D1
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2 data: C2 becomes the species text (previously duplicate numeric max),
# D2 stays the same species text, E2 becomes the numeric rejection-f value
$ws.Range("C2").Value = "s__CAG-603 sp900066105"
$ws.Range("D2").Value = "s__CAG-603 sp900066105"
$ws.Range("E2").Value = 0.9999999999991507
